# Slow test refactor: Ignore Excel tests > 100ms
#
# Applies the data edits captured in the target diff:
#  - TestRecord!A10  43218 -> 43220
#  - TestRecord!B10  67.14 -> 69.540000000000006
#  - TestRecord!E10  "some test textzzz...(39 z's)" -> "some test textzzz...(41 z's)"
#  - Budget Out!C9   83.02 -> 84.22
#  - Budget Out!F9   "Description007zzz...(38 z's)" -> "Description007zzz...(39 z's)"
#  - Expected Out!B9  1324.16 -> 1326.56
#  - Expected Out!B11 420.82  -> 422.02
#  - Expected Out!B1 (=SUM(B2:B295)) recalculates automatically from the above
#  - Window position (cosmetic, best effort)

$wb = $excel.ActiveWorkbook

# Cosmetic: workbook window position moved on last save.
$win = $excel.ActiveWindow
$win.Left = 780
$win.Top = 780

# --- TestRecord sheet ---
$testRecord = $wb.Worksheets.Item("TestRecord")
$testRecord.Range("A10").Value = 43220
$testRecord.Range("B10").Value = 69.540000000000006
$testRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Budget Out sheet ---
$budgetOut = $wb.Worksheets.Item("Budget Out")
$budgetOut.Range("C9").Value = 84.22
$budgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Expected Out sheet ---
$expectedOut = $wb.Worksheets.Item("Expected Out")
$expectedOut.Range("B9").Value = 1326.56
$expectedOut.Range("B11").Value = 422.02

$excel.CalculateFull()

Write-Host "Edits applied"
